# nodes.xlsx — "added enekos dv, refactored"
# - D4 was miscategorised as "hospital"; it is actually a "base".
# - E5/E6 priorities reset to 0 (no longer the only/last hospitals).
# - Five new "hospital" nodes (the Enekos batch) are appended as rows 7-11,
#   continuing the A-column running index via the same fill-down formula
#   used by the existing rows, with B/C lat-long and E priority values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 reclassified from hospital -> base
$ws.Range("D4").Value = "base"

# Row 5/6 priority reset to 0
$ws.Range("E5").Value = 0
$ws.Range("E6").Value = 0

# New row 7 (Enekos)
$ws.Range("A7").Formula = "=A6+1"
$ws.Range("B7").Value = 6.299418
$ws.Range("C7").Value = 81.054685
$ws.Range("D7").Value = "hospital"
$ws.Range("E7").Value = 1

# New row 8
$ws.Range("A8").Formula = "=A7+1"
$ws.Range("B8").Value = 6.496502
$ws.Range("C8").Value = 80.776137
$ws.Range("D8").Value = "hospital"
$ws.Range("E8").Value = 1

# New row 9
$ws.Range("A9").Formula = "=A8+1"
$ws.Range("B9").Value = 6.755147
$ws.Range("C9").Value = 80.971558
$ws.Range("D9").Value = "hospital"
$ws.Range("E9").Value = 1

# New row 10
$ws.Range("A10").Formula = "=A9+1"
$ws.Range("B10").Value = 6.978419
$ws.Range("C10").Value = 81.246235
$ws.Range("D10").Value = "hospital"
$ws.Range("E10").Value = 2

# New row 11
$ws.Range("A11").Formula = "=A10+1"
$ws.Range("B11").Value = 7.420607
$ws.Range("C11").Value = 81.590763
$ws.Range("D11").Value = "hospital"
$ws.Range("E11").Value = 2

# Cursor ends up on G12, matching the saved selection in the workbook
$ws.Range("G12").Select()
